# Generate Report for Handoff
# Adds a new file (e577379d-7c22-421d-8b27-1a321313dc7b) as a new row on each
# of the three sheets (Overview, zh-cn, de-de), mirroring the existing row
# for 7e859f56-bd5d-495d-bf75-ac5f3cad680c.

$wb = $excel.ActiveWorkbook

$newGuid = "e577379d-7c22-421d-8b27-1a321313dc7b"
$newHash = "53df4f705f2a771179cfdce892b56b72e6409af5"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "$newGuid.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-28-19 12:28:00"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/74d0d9ee8cd7cab6423a8350304d58897262a311/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "$newGuid.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-19 12:27:57"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/74d0d9ee8cd7cab6423a8350304d58897262a311/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/74d0d9ee8cd7cab6423a8350304d58897262a311/e2e/$newGuid.md",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f0f3aad4f16d91a281c7c4cc9b54e499ef55d139/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf",
    "",
    "",
    "$newGuid.$newHash.zh-cn.xlf"
)

$wsZhCn.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "$newGuid.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-19 12:28:00"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/74d0d9ee8cd7cab6423a8350304d58897262a311/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/74d0d9ee8cd7cab6423a8350304d58897262a311/e2e/$newGuid.md",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3e8214e59b2fdd0395b808413b75797a3fad1a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf",
    "",
    "",
    "$newGuid.$newHash.de-de.xlf"
)

$wsDeDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Output "Applied handoff report update for $newGuid"
